$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Group"
$ws.Range("B2").Value = "Jul"
$ws.Range("E2").Value = 28.4821
$ws.Range("F2").Value = 0.0
$ws.Range("B3").Value = "Aug"
$ws.Range("E3").Value = 35.4163
$ws.Range("F3").Value = 0.0
$ws.Range("B4").Value = "Sep"
$ws.Range("E4").Value = 20.5463
$ws.Range("F4").Value = 0.0
$ws.Range("B5").Value = "Oct"
$ws.Range("E5").Value = 55.1314
$ws.Range("F5").Value = 0.0
$ws.Range("B6").Value = "Nov"
$ws.Range("E6").Value = 47.5583
$ws.Range("F6").Value = 0.0
$ws.Range("B7").Value = "Dec"
$ws.Range("E7").Value = 65.1021
$ws.Range("F7").Value = 0.0
$ws.Range("B8").Value = "Jan"
$ws.Range("E8").Value = 71.6884
$ws.Range("F8").Value = 0.0
$ws.Range("B9").Value = "Feb"
$ws.Range("E9").Value = 91.3516
$ws.Range("F9").Value = 0.0
$ws.Range("B10").Value = "Mar"
$ws.Range("E10").Value = 139.722
$ws.Range("F10").Value = 0.0
$ws.Range("B11").Value = "Apr"
$ws.Range("E11").Value = 81.8066
$ws.Range("F11").Value = 0.0
$ws.Range("B12").Value = "May"
$ws.Range("E12").Value = 51.8385
$ws.Range("F12").Value = 0.0
$ws.Range("B13").Value = "Jun"
$ws.Range("E13").Value = 4.1695
$ws.Range("F13").Value = 0.1243
$ws.Range("E14").Value = 643.6463
$ws.Range("F14").Value = 0.0
$ws.Range("B15").Value = "Jul vs. Other Months"
$ws.Range("E15").Value = 0.0891
$ws.Range("F15").Value = 0.7654
$ws.Range("B16").Value = "Aug vs. Other Months"
$ws.Range("E16").Value = 2.3118
$ws.Range("F16").Value = 0.1288
$ws.Range("B17").Value = "Sep vs. Other Months"
$ws.Range("E17").Value = 0.4164
$ws.Range("F17").Value = 0.5189
$ws.Range("B18").Value = "Oct vs. Other Months"
$ws.Range("E18").Value = 0.0495
$ws.Range("F18").Value = 0.8239
$ws.Range("B19").Value = "Nov vs. Other Months"
$ws.Range("E19").Value = 0.0011
$ws.Range("F19").Value = 0.9735
$ws.Range("B20").Value = "Dec vs. Other Months"
$ws.Range("E20").Value = 0.1567
$ws.Range("F20").Value = 0.6923
$ws.Range("B21").Value = "Jan vs. Other Months"
$ws.Range("E21").Value = 0.8724
$ws.Range("F21").Value = 0.3506
$ws.Range("B22").Value = "Feb vs. Other Months"
$ws.Range("E22").Value = 1.3965
$ws.Range("F22").Value = 0.2376
$ws.Range("B23").Value = "Mar vs. Other Months"
$ws.Range("E23").Value = 0.5157
$ws.Range("F23").Value = 0.4729
$ws.Range("B24").Value = "Apr vs. Other Months"
$ws.Range("E24").Value = 1.6496
$ws.Range("F24").Value = 0.1994
$ws.Range("B25").Value = "May vs. Other Months"
$ws.Range("E25").Value = 0.0238
$ws.Range("F25").Value = 0.8773
$ws.Range("B26").Value = "Jun vs. Other Months"
$ws.Range("E26").Value = 1.9309
$ws.Range("F26").Value = 0.165
$ws.Range("B27").Value = "Jul vs. Other Months"
$ws.Range("F27").Value = 0.0211
$ws.Range("B28").Value = "Aug vs. Other Months"
$ws.Range("F28").Value = 0.0168
$ws.Range("B29").Value = "Sep vs. Other Months"
$ws.Range("F29").Value = 0.0125
$ws.Range("B30").Value = "Oct vs. Other Months"
$ws.Range("F30").Value = 0.3573
$ws.Range("B31").Value = "Nov vs. Other Months"
$ws.Range("F31").Value = 0.0754
$ws.Range("B32").Value = "Dec vs. Other Months"
$ws.Range("F32").Value = 0.107
$ws.Range("B33").Value = "Jan vs. Other Months"
$ws.Range("F33").Value = 0.2502
$ws.Range("B34").Value = "Feb vs. Other Months"
$ws.Range("F34").Value = 0.0011
$ws.Range("B35").Value = "Mar vs. Other Months"
$ws.Range("F35").Value = 0.0001
$ws.Range("B36").Value = "Apr vs. Other Months"
$ws.Range("F36").Value = 0.0137
$ws.Range("B37").Value = "May vs. Other Months"
$ws.Range("F37").Value = 0.1652
$ws.Range("B38").Value = "Jun vs. Other Months"
$ws.Range("F38").Value = 0.1355
$ws.Range("B39").Value = "Jul"
$ws.Range("B40").Value = "Aug"
$ws.Range("B41").Value = "Sep"
$ws.Range("B42").Value = "Oct"
$ws.Range("B43").Value = "Nov"
$ws.Range("B44").Value = "Dec"
$ws.Range("B45").Value = "Jan"
$ws.Range("B46").Value = "Feb"
$ws.Range("B47").Value = "Mar"
$ws.Range("B48").Value = "Apr"
$ws.Range("B49").Value = "May"
$ws.Range("B50").Value = "Jun"
